$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("Sheet2")
$ws3 = $wb.Worksheets.Item("Sheet3")

# Preserve Sheet3's leftover selection (full columns O:P) before we move away from it,
# matching the state captured in the target workbook.
$ws3.Activate()
$ws3.Range("O1:P1048576").Select()

# Switch to Sheet2 - this is where the new "drop" / "dropExplanation" columns are added.
$ws2.Activate()

# New header cells Q1 / R1, copying the header formatting already used by column P.
$ws2.Range("P1").Copy()
$ws2.Range("Q1:R1").PasteSpecial(-4122)

$ws2.Range("Q1").Value = "drop"
$ws2.Range("R1").Value = "dropExplanation"

# Every data row (2-41) gets a FALSE boolean value in the new "drop" column.
for ($r = 2; $r -le 41; $r++) {
    $ws2.Cells.Item($r, 17).Value = $false
}

# Leave the full columns Q:R selected, as in the source edit.
$ws2.Range("Q1:R1048576").Select()
